# Generate Report for Handback
$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text "Ready for handoff" -> "Handed back: in sync with en-US"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: add Latest Target File / Latest Handback File columns (F, G),
#     update Status (C) and Latest Handback DateTime (H)
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/194f448d87c97ac1354b4a4819c98310c4a62bb5/e2e/22553294-dc54-45e9-a035-b9ff982c7246.md", "", "", "22553294-dc54-45e9-a035-b9ff982c7246.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1439a17a12f3e9e1a71174e5c93a5731045c5bdf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/22553294-dc54-45e9-a035-b9ff982c7246.02b9fb26a72705b3c083373750eb05a67a9d3c13.zh-cn.xlf", "", "", "22553294-dc54-45e9-a035-b9ff982c7246.02b9fb26a72705b3c083373750eb05a67a9d3c13.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/194f448d87c97ac1354b4a4819c98310c4a62bb5/e2e/83b3158e-289f-43c5-b85d-30d454e7233b.md", "", "", "83b3158e-289f-43c5-b85d-30d454e7233b.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1439a17a12f3e9e1a71174e5c93a5731045c5bdf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/83b3158e-289f-43c5-b85d-30d454e7233b.18584ae04e2afaa310860b5ee7fa12982a2be895.zh-cn.xlf", "", "", "83b3158e-289f-43c5-b85d-30d454e7233b.18584ae04e2afaa310860b5ee7fa12982a2be895.zh-cn.xlf")

$wsZhCn.Range("H2").Value = "2016-03-12 06:42:22"
$wsZhCn.Range("H3").Value = "2016-03-12 06:42:22"

# --- de-de sheet: add Latest Target File / Latest Handback File columns (F, G),
#     update Status (C) and Latest Handback DateTime (H)
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/194f448d87c97ac1354b4a4819c98310c4a62bb5/e2e/22553294-dc54-45e9-a035-b9ff982c7246.md", "", "", "22553294-dc54-45e9-a035-b9ff982c7246.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a723dc40bb5609f761fb95755f980a2f6d89a56d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/22553294-dc54-45e9-a035-b9ff982c7246.02b9fb26a72705b3c083373750eb05a67a9d3c13.de-de.xlf", "", "", "22553294-dc54-45e9-a035-b9ff982c7246.02b9fb26a72705b3c083373750eb05a67a9d3c13.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/194f448d87c97ac1354b4a4819c98310c4a62bb5/e2e/83b3158e-289f-43c5-b85d-30d454e7233b.md", "", "", "83b3158e-289f-43c5-b85d-30d454e7233b.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a723dc40bb5609f761fb95755f980a2f6d89a56d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/83b3158e-289f-43c5-b85d-30d454e7233b.18584ae04e2afaa310860b5ee7fa12982a2be895.de-de.xlf", "", "", "83b3158e-289f-43c5-b85d-30d454e7233b.18584ae04e2afaa310860b5ee7fa12982a2be895.de-de.xlf")

$wsDeDe.Range("H2").Value = "2016-03-12 06:42:28"
$wsDeDe.Range("H3").Value = "2016-03-12 06:42:28"
